$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overall")

$ws.Range("H26").Value = 0.74326
$ws.Range("I26").Value = 0.02266
$ws.Range("H27").Value = 0.49204
$ws.Range("I27").Value = 0.05956
$ws.Range("H28").Value = 0.74416
$ws.Range("I28").Value = 0.02283
$ws.Range("H29").Value = 0.49588
$ws.Range("I29").Value = 0.0601
$ws.Range("H30").Value = 0.74307
$ws.Range("I30").Value = 0.02306
$ws.Range("H31").Value = 0.48934
$ws.Range("I31").Value = 0.06162
$ws.Range("H32").Value = 0.74468
$ws.Range("I32").Value = 0.02333
$ws.Range("H33").Value = 0.49588
$ws.Range("I33").Value = 0.0609
$ws.Range("H34").Value = 0.72419
$ws.Range("I34").Value = 0.02342
$ws.Range("H35").Value = 0.38992
$ws.Range("I35").Value = 0.06596
$ws.Range("H36").Value = 0.72471
$ws.Range("I36").Value = 0.02328
$ws.Range("H37").Value = 0.393
$ws.Range("I37").Value = 0.06646000000000001
$ws.Range("H38").Value = 0.724
$ws.Range("I38").Value = 0.02339
$ws.Range("H39").Value = 0.38992
$ws.Range("I39").Value = 0.06596
$ws.Range("H40").Value = 0.72458
$ws.Range("I40").Value = 0.02325
$ws.Range("H41").Value = 0.393
$ws.Range("I41").Value = 0.06646000000000001
$ws.Range("H66").Value = 0.7212
$ws.Range("I66").Value = 0.02794
$ws.Range("H67").Value = 0.42037
$ws.Range("I67").Value = 0.05307
$ws.Range("H68").Value = 0.7309099999999999
$ws.Range("I68").Value = 0.02772
$ws.Range("H69").Value = 0.47586
$ws.Range("I69").Value = 0.05948
$ws.Range("H70").Value = 0.72088
$ws.Range("I70").Value = 0.02801
$ws.Range("H71").Value = 0.41959
$ws.Range("I71").Value = 0.05376
$ws.Range("H72").Value = 0.73117
$ws.Range("I72").Value = 0.0284
$ws.Range("H73").Value = 0.48087
$ws.Range("I73").Value = 0.05803
$ws.Range("H74").Value = 0.70939
$ws.Range("I74").Value = 0.02667
$ws.Range("H75").Value = 0.35948
$ws.Range("I75").Value = 0.04855
$ws.Range("H76").Value = 0.71517
$ws.Range("I76").Value = 0.02772
$ws.Range("H77").Value = 0.39299
$ws.Range("I77").Value = 0.05445
$ws.Range("H78").Value = 0.70952
$ws.Range("I78").Value = 0.02636
$ws.Range("H79").Value = 0.35948
$ws.Range("I79").Value = 0.04855
$ws.Range("H80").Value = 0.71569
$ws.Range("I80").Value = 0.02717
$ws.Range("H81").Value = 0.39415
$ws.Range("I81").Value = 0.05424
$ws.Range("H106").Value = 0.7401
$ws.Range("I106").Value = 0.01886
$ws.Range("H107").Value = 0.47475
$ws.Range("I107").Value = 0.05502
$ws.Range("H108").Value = 0.74266
$ws.Range("I108").Value = 0.02103
$ws.Range("H109").Value = 0.48477
$ws.Range("I109").Value = 0.05843
$ws.Range("H110").Value = 0.73991
$ws.Range("I110").Value = 0.01962
$ws.Range("H111").Value = 0.47514
$ws.Range("I111").Value = 0.05629
$ws.Range("H112").Value = 0.74234
$ws.Range("I112").Value = 0.02082
$ws.Range("H113").Value = 0.48438
$ws.Range("I113").Value = 0.06027
$ws.Range("H114").Value = 0.71956
$ws.Range("I114").Value = 0.01894
$ws.Range("H115").Value = 0.35799
$ws.Range("I115").Value = 0.06148
$ws.Range("H116").Value = 0.71943
$ws.Range("I116").Value = 0.01934
$ws.Range("H117").Value = 0.36069
$ws.Range("I117").Value = 0.0614
$ws.Range("H118").Value = 0.71962
$ws.Range("I118").Value = 0.01893
$ws.Range("H119").Value = 0.35799
$ws.Range("I119").Value = 0.06148
$ws.Range("H120").Value = 0.71962
$ws.Range("I120").Value = 0.0194
$ws.Range("H121").Value = 0.36069
$ws.Range("I121").Value = 0.0614
$ws.Range("H146").Value = 0.73774
$ws.Range("I146").Value = 0.02739
$ws.Range("H147").Value = 0.46469
$ws.Range("I147").Value = 0.06397
$ws.Range("H148").Value = 0.74418
$ws.Range("I148").Value = 0.02759
$ws.Range("H149").Value = 0.50479
$ws.Range("I149").Value = 0.06582
$ws.Range("H150").Value = 0.73793
$ws.Range("I150").Value = 0.02756
$ws.Range("H151").Value = 0.46469
$ws.Range("I151").Value = 0.06397
$ws.Range("H152").Value = 0.7474499999999999
$ws.Range("I152").Value = 0.02656
$ws.Range("H153").Value = 0.51019
$ws.Range("I153").Value = 0.06537999999999999
$ws.Range("H154").Value = 0.72104
$ws.Range("I154").Value = 0.02549
$ws.Range("H155").Value = 0.38228
$ws.Range("I155").Value = 0.05566
$ws.Range("H156").Value = 0.72412
$ws.Range("I156").Value = 0.02616
$ws.Range("H157").Value = 0.4027
$ws.Range("I157").Value = 0.05478
$ws.Range("H158").Value = 0.72168
$ws.Range("I158").Value = 0.02515
$ws.Range("H159").Value = 0.38189
$ws.Range("I159").Value = 0.05528
$ws.Range("H160").Value = 0.72502
$ws.Range("I160").Value = 0.0268
$ws.Range("H161").Value = 0.40462
$ws.Range("I161").Value = 0.05393
